# Update "想去人数" (want-to-go count) figures for several rows across
# the 展览 (Exhibition), 本地生活 (Local life), and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# 展览 sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 1138
$ws1.Range("F9").Value = 228
$ws1.Range("F11").Value = 8004
$ws1.Range("F13").Value = 9497
$ws1.Range("F29").Value = 1625
$ws1.Range("F34").Value = 42
$ws1.Range("F46").Value = 41
$ws1.Range("F48").Value = 90

# 本地生活 sheet
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2730

# 全部类型 sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 1138
$ws4.Range("F17").Value = 8004
$ws4.Range("F18").Value = 9497
$ws4.Range("F28").Value = 1625
$ws4.Range("F45").Value = 41
$ws4.Range("F49").Value = 90
